$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: the phone number for payment 71277620 settles from a text value
# into a plain numeric value (A33).
$ws.Cells.Item(33, 1).Value = 71277620

# Row 34: append the new payment 71277620 (Cash) 2025-08-18T17:10:08
$ws.Cells.Item(34, 1).NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "71277620"
$ws.Cells.Item(34, 1).Style = "Normal"

$ws.Cells.Item(34, 3).Value = "Cash"
$ws.Cells.Item(34, 4).Value = "2025-08-18T17:10:08"
$ws.Cells.Item(34, 5).Value = 76
$ws.Cells.Item(34, 7).Value = 76
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
